# Update page count for "Title pages" from 2 to 4.
# This reflects the extra page break added to the title page so that the
# Title page and Table of Contents now appear on right hand pages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4

# Update the active cell selection on the sheet (bottom-left pane of the
# frozen view) to C3, matching the author's cursor position when saving.
$ws.Range("C3").Select()
